$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: insert a new "In Bearbeitung" (In Progress) column header in
#    D2, and move the existing "Erledigt" (Done) header text over to E2.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "In Bearbeitung"
$ws.Range("E2").Value = "Erledigt"

# ---------------------------------------------------------------------------
# 2. Give the new "In Bearbeitung" status column (D3:D10) its own visual
#    style: 14pt font, centered horizontally/vertically, wrapped text.
#    Apply it to a single cell first (keeps the style table clean), then
#    copy/paste the format onto the rest of the column so every cell ends
#    up sharing the very same cell-format entry instead of Excel fanning
#    out a new xf per cell.
# ---------------------------------------------------------------------------
$ws.Range("D3").Font.Size = 14
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").VerticalAlignment = -4108
$ws.Range("D3").WrapText = $true

$ws.Range("D3").Copy()
$ws.Range("D4:D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Data updates for the new column / assignment changes.
# ---------------------------------------------------------------------------
# Mark row 7 ("Shader Programmierung") as currently in progress.
$ws.Range("D7").Value = "X"

# Assign Joe as the second Bearbeiter (person) on row 8 ("Aktionen des
# Players Fkt.").
$ws.Range("C8").Value = "Joe"

# ---------------------------------------------------------------------------
# 4. Column / row sizing.
# ---------------------------------------------------------------------------
# Give column D (the new status column) a custom width.
$ws.Columns.Item(4).ColumnWidth = 13.9

# Row 4 gets a slightly taller custom height.
$ws.Rows.Item(4).RowHeight = 18.75

# ---------------------------------------------------------------------------
# 5. Selection moves to G10.
# ---------------------------------------------------------------------------
$ws.Range("G10").Select()
